$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph 1 ("This is a Microsoft word document.") gets two
# trailing spaces appended to its existing run, followed by three new runs
# (all colored C00000 / dark red) spelling out the "(This is a change - ...)"
# annotation.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">This is a Microsoft word document.  </w:t></w:r><w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>(This is a change – Ve</w:t></w:r><w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>rsion for branch alternate</w:t></w:r><w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: the blank paragraph between "It will be treated..." and
# "The Raven" gains paragraph-mark formatting (light-grey shading, Calibri,
# bold, dark-grey text color) - it stays empty, just restyled.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="202122"/></w:rPr></w:pPr></w:p>
'@
$p3.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Change 3: the trailing "ank God almighty, we are free at last." paragraph
# (NormalWeb style) is cleared out to a plain empty paragraph.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$pLast = $d.Paragraphs($n)
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$pLast.Range.InsertXML($xml3)

Write-Output "edits applied"
